# Rename GenericPlate cargo names: drop the "_idN" plate-id segment from
# antiNelson/antiQuimby entries and fold the "_y0n9m4" suffix into a
# hyphenated name prefix for loopctrl/mask entries. Also update the
# header cell (A1) that documents the naming pattern, and move the
# on-sheet selection to where the author last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Names")
$ws.Activate()

# A1 documents the naming scheme used by the columns below.
$ws.Range("A1").Value = "name-side-position"

# Cargo / control / mask handle names for h2, rows 2-9 (columns B:Y).
$ws.Range("B2").Value = "antiNelson_h2_pos1"
$ws.Range("C2").Value = "antiNelson_h2_pos2"
$ws.Range("D2").Value = "antiNelson_h2_pos3"
$ws.Range("E2").Value = "antiNelson_h2_pos4"
$ws.Range("F2").Value = "antiNelson_h2_pos5"
$ws.Range("G2").Value = "antiNelson_h2_pos6"
$ws.Range("H2").Value = "antiNelson_h2_pos7"
$ws.Range("I2").Value = "antiNelson_h2_pos8"
$ws.Range("J2").Value = "antiNelson_h2_pos9"
$ws.Range("K2").Value = "antiNelson_h2_pos10"
$ws.Range("L2").Value = "antiNelson_h2_pos11"
$ws.Range("M2").Value = "antiNelson_h2_pos12"
$ws.Range("N2").Value = "antiNelson_h2_pos13"
$ws.Range("O2").Value = "antiNelson_h2_pos14"
$ws.Range("P2").Value = "antiNelson_h2_pos15"
$ws.Range("Q2").Value = "antiNelson_h2_pos16"
$ws.Range("R2").Value = "antiNelson_h2_pos17"
$ws.Range("S2").Value = "antiNelson_h2_pos18"
$ws.Range("T2").Value = "antiNelson_h2_pos19"
$ws.Range("U2").Value = "antiNelson_h2_pos20"
$ws.Range("V2").Value = "antiNelson_h2_pos21"
$ws.Range("W2").Value = "antiNelson_h2_pos22"
$ws.Range("X2").Value = "antiNelson_h2_pos23"
$ws.Range("Y2").Value = "antiNelson_h2_pos24"
$ws.Range("B3").Value = "antiNelson_h2_pos25"
$ws.Range("C3").Value = "antiNelson_h2_pos26"
$ws.Range("D3").Value = "antiNelson_h2_pos27"
$ws.Range("E3").Value = "antiNelson_h2_pos28"
$ws.Range("F3").Value = "antiNelson_h2_pos29"
$ws.Range("G3").Value = "antiNelson_h2_pos30"
$ws.Range("H3").Value = "antiNelson_h2_pos31"
$ws.Range("I3").Value = "antiNelson_h2_pos32"
$ws.Range("B4").Value = "antiQuimby_h2_pos1"
$ws.Range("C4").Value = "antiQuimby_h2_pos2"
$ws.Range("D4").Value = "antiQuimby_h2_pos3"
$ws.Range("E4").Value = "antiQuimby_h2_pos4"
$ws.Range("F4").Value = "antiQuimby_h2_pos5"
$ws.Range("G4").Value = "antiQuimby_h2_pos6"
$ws.Range("H4").Value = "antiQuimby_h2_pos7"
$ws.Range("I4").Value = "antiQuimby_h2_pos8"
$ws.Range("J4").Value = "antiQuimby_h2_pos9"
$ws.Range("K4").Value = "antiQuimby_h2_pos10"
$ws.Range("L4").Value = "antiQuimby_h2_pos11"
$ws.Range("M4").Value = "antiQuimby_h2_pos12"
$ws.Range("N4").Value = "antiQuimby_h2_pos13"
$ws.Range("O4").Value = "antiQuimby_h2_pos14"
$ws.Range("P4").Value = "antiQuimby_h2_pos15"
$ws.Range("Q4").Value = "antiQuimby_h2_pos16"
$ws.Range("R4").Value = "antiQuimby_h2_pos17"
$ws.Range("S4").Value = "antiQuimby_h2_pos18"
$ws.Range("T4").Value = "antiQuimby_h2_pos19"
$ws.Range("U4").Value = "antiQuimby_h2_pos20"
$ws.Range("V4").Value = "antiQuimby_h2_pos21"
$ws.Range("W4").Value = "antiQuimby_h2_pos22"
$ws.Range("X4").Value = "antiQuimby_h2_pos23"
$ws.Range("Y4").Value = "antiQuimby_h2_pos24"
$ws.Range("B5").Value = "antiQuimby_h2_pos25"
$ws.Range("C5").Value = "antiQuimby_h2_pos26"
$ws.Range("D5").Value = "antiQuimby_h2_pos27"
$ws.Range("E5").Value = "antiQuimby_h2_pos28"
$ws.Range("F5").Value = "antiQuimby_h2_pos29"
$ws.Range("G5").Value = "antiQuimby_h2_pos30"
$ws.Range("H5").Value = "antiQuimby_h2_pos31"
$ws.Range("I5").Value = "antiQuimby_h2_pos32"
$ws.Range("B6").Value = "loopctrl-y0n9m4_h2_pos1"
$ws.Range("C6").Value = "loopctrl-y0n9m4_h2_pos2"
$ws.Range("D6").Value = "loopctrl-y0n9m4_h2_pos3"
$ws.Range("E6").Value = "loopctrl-y0n9m4_h2_pos4"
$ws.Range("F6").Value = "loopctrl-y0n9m4_h2_pos5"
$ws.Range("G6").Value = "loopctrl-y0n9m4_h2_pos6"
$ws.Range("H6").Value = "loopctrl-y0n9m4_h2_pos7"
$ws.Range("I6").Value = "loopctrl-y0n9m4_h2_pos8"
$ws.Range("J6").Value = "loopctrl-y0n9m4_h2_pos9"
$ws.Range("K6").Value = "loopctrl-y0n9m4_h2_pos10"
$ws.Range("L6").Value = "loopctrl-y0n9m4_h2_pos11"
$ws.Range("M6").Value = "loopctrl-y0n9m4_h2_pos12"
$ws.Range("N6").Value = "loopctrl-y0n9m4_h2_pos13"
$ws.Range("O6").Value = "loopctrl-y0n9m4_h2_pos14"
$ws.Range("P6").Value = "loopctrl-y0n9m4_h2_pos15"
$ws.Range("Q6").Value = "loopctrl-y0n9m4_h2_pos16"
$ws.Range("B7").Value = "loopctrl-y0n9m4_h2_pos17"
$ws.Range("C7").Value = "loopctrl-y0n9m4_h2_pos18"
$ws.Range("D7").Value = "loopctrl-y0n9m4_h2_pos19"
$ws.Range("E7").Value = "loopctrl-y0n9m4_h2_pos20"
$ws.Range("F7").Value = "loopctrl-y0n9m4_h2_pos21"
$ws.Range("G7").Value = "loopctrl-y0n9m4_h2_pos22"
$ws.Range("H7").Value = "loopctrl-y0n9m4_h2_pos23"
$ws.Range("I7").Value = "loopctrl-y0n9m4_h2_pos24"
$ws.Range("J7").Value = "loopctrl-y0n9m4_h2_pos25"
$ws.Range("K7").Value = "loopctrl-y0n9m4_h2_pos26"
$ws.Range("L7").Value = "loopctrl-y0n9m4_h2_pos27"
$ws.Range("M7").Value = "loopctrl-y0n9m4_h2_pos28"
$ws.Range("N7").Value = "loopctrl-y0n9m4_h2_pos29"
$ws.Range("O7").Value = "loopctrl-y0n9m4_h2_pos30"
$ws.Range("P7").Value = "loopctrl-y0n9m4_h2_pos31"
$ws.Range("Q7").Value = "loopctrl-y0n9m4_h2_pos32"
$ws.Range("B8").Value = "mask-y0n9m4_h2_pos1"
$ws.Range("C8").Value = "mask-y0n9m4_h2_pos2"
$ws.Range("D8").Value = "mask-y0n9m4_h2_pos3"
$ws.Range("E8").Value = "mask-y0n9m4_h2_pos4"
$ws.Range("F8").Value = "mask-y0n9m4_h2_pos5"
$ws.Range("G8").Value = "mask-y0n9m4_h2_pos6"
$ws.Range("H8").Value = "mask-y0n9m4_h2_pos7"
$ws.Range("I8").Value = "mask-y0n9m4_h2_pos8"
$ws.Range("J8").Value = "mask-y0n9m4_h2_pos9"
$ws.Range("K8").Value = "mask-y0n9m4_h2_pos10"
$ws.Range("L8").Value = "mask-y0n9m4_h2_pos11"
$ws.Range("M8").Value = "mask-y0n9m4_h2_pos12"
$ws.Range("N8").Value = "mask-y0n9m4_h2_pos13"
$ws.Range("O8").Value = "mask-y0n9m4_h2_pos14"
$ws.Range("P8").Value = "mask-y0n9m4_h2_pos15"
$ws.Range("Q8").Value = "mask-y0n9m4_h2_pos16"
$ws.Range("B9").Value = "mask-y0n9m4_h2_pos17"
$ws.Range("C9").Value = "mask-y0n9m4_h2_pos18"
$ws.Range("D9").Value = "mask-y0n9m4_h2_pos19"
$ws.Range("E9").Value = "mask-y0n9m4_h2_pos20"
$ws.Range("F9").Value = "mask-y0n9m4_h2_pos21"
$ws.Range("G9").Value = "mask-y0n9m4_h2_pos22"
$ws.Range("H9").Value = "mask-y0n9m4_h2_pos23"
$ws.Range("I9").Value = "mask-y0n9m4_h2_pos24"
$ws.Range("J9").Value = "mask-y0n9m4_h2_pos25"
$ws.Range("K9").Value = "mask-y0n9m4_h2_pos26"
$ws.Range("L9").Value = "mask-y0n9m4_h2_pos27"
$ws.Range("M9").Value = "mask-y0n9m4_h2_pos28"
$ws.Range("N9").Value = "mask-y0n9m4_h2_pos29"
$ws.Range("O9").Value = "mask-y0n9m4_h2_pos30"
$ws.Range("P9").Value = "mask-y0n9m4_h2_pos31"
$ws.Range("Q9").Value = "mask-y0n9m4_h2_pos32"

# Leave the selection where the author left it after editing.
$ws.Range("O27").Select()
